$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting old row 6 (Ensemble) down to row 7
$ws.Rows(6).Insert()

# Set label and style for new RF row (A6) to match other label cells
$ws.Range("A6").Value = "RF"
$a6 = $ws.Range("A6")
$a6.Font.Bold = $true
$a6.HorizontalAlignment = -4108
$a6.VerticalAlignment = -4160
$a6.Borders.LineStyle = 1
$a6.Borders.Weight = 2

# Row 2 values
$ws.Range("B2").Value = 0.6109102501821633
$ws.Range("C2").Value = 0.05756199034196108
$ws.Range("D2").Value = 0.5712892319620031
$ws.Range("E2").Value = 0.5914397993824438
$ws.Range("F2").Value = 0.5852399816746165
$ws.Range("G2").Value = 0.5812954546253432
$ws.Range("H2").Value = 0.7252867832664098
$ws.Range("I2").Value = 0.677086542942525
$ws.Range("J2").Value = 0.02873757055082129
$ws.Range("K2").Value = 0.6645249064177898
$ws.Range("L2").Value = 0.6338127682213704
$ws.Range("M2").Value = 0.6922143266360901
$ws.Range("N2").Value = 0.6746212835005254
$ws.Range("O2").Value = 0.7202594299368493
$ws.Range("P2").Value = 0.6754991818369545
$ws.Range("Q2").Value = 0.05157617898724755
$ws.Range("R2").Value = 0.6958576927991821
$ws.Range("S2").Value = 0.6092563056730569
$ws.Range("T2").Value = 0.7354036063713483
$ws.Range("U2").Value = 0.6194427329694168
$ws.Range("V2").Value = 0.7175355713717684
$ws.Range("W2").Value = 0.791567739767727
$ws.Range("X2").Value = 0.03740029654046373
$ws.Range("Y2").Value = 0.8194955649548281
$ws.Range("Z2").Value = 0.7271344102551904
$ws.Range("AA2").Value = 0.782035349224515
$ws.Range("AB2").Value = 0.8359371343095949
$ws.Range("AC2").Value = 0.7932362400945066
$ws.Range("AD2").Value = 0.7108778336318742
$ws.Range("AE2").Value = 0.03707655289613247
$ws.Range("AF2").Value = 0.6545585673233221
$ws.Range("AG2").Value = 0.729531490015361
$ws.Range("AH2").Value = 0.7014149064869292
$ws.Range("AI2").Value = 0.7672555043522785
$ws.Range("AJ2").Value = 0.7016286999814797
$ws.Range("AK2").Value = 0.6977995982724715
$ws.Range("AL2").Value = 0.04350338369607885
$ws.Range("AM2").Value = 0.6378493505432825
$ws.Range("AN2").Value = 0.7098914285178592
$ws.Range("AO2").Value = 0.6669494560070145
$ws.Range("AP2").Value = 0.7660781537375686
$ws.Range("AQ2").Value = 0.7082296025566325

# Row 3 values
$ws.Range("B3").Value = 0.8008134154112974
$ws.Range("C3").Value = 0.02802671679896426
$ws.Range("D3").Value = 0.7996874646398039
$ws.Range("E3").Value = 0.7746018616986359
$ws.Range("F3").Value = 0.8072219710693352
$ws.Range("G3").Value = 0.7726832598587123
$ws.Range("H3").Value = 0.8498725197899992
$ws.Range("I3").Value = 0.8089305373102448
$ws.Range("J3").Value = 0.04992339510573657
$ws.Range("K3").Value = 0.7487308543727752
$ws.Range("L3").Value = 0.7691540769726628
$ws.Range("M3").Value = 0.8301849850236948
$ws.Range("N3").Value = 0.8052492310556826
$ws.Range("O3").Value = 0.8913335391264083
$ws.Range("P3").Value = 0.7993243619866995
$ws.Range("Q3").Value = 0.05106254769302052
$ws.Range("R3").Value = 0.788865582459477
$ws.Range("S3").Value = 0.7239145920934016
$ws.Range("T3").Value = 0.8389012969213979
$ws.Range("U3").Value = 0.7744956507053281
$ws.Range("V3").Value = 0.8704446877538931
$ws.Range("W3").Value = 0.8307171212580075
$ws.Range("X3").Value = 0.04124691612600012
$ws.Range("Y3").Value = 0.8627762972558537
$ws.Range("Z3").Value = 0.75426020247752
$ws.Range("AA3").Value = 0.8387134088322543
$ws.Range("AB3").Value = 0.8276268217751301
$ws.Range("AC3").Value = 0.8702088759492792
$ws.Range("AD3").Value = 0.8475142028388465
$ws.Range("AE3").Value = 0.03345178412042464
$ws.Range("AF3").Value = 0.8511172266491415
$ws.Range("AG3").Value = 0.7852325310429774
$ws.Range("AH3").Value = 0.8709406499729081
$ws.Range("AI3").Value = 0.84881603261981
$ws.Range("AJ3").Value = 0.8814645739093956
$ws.Range("AK3").Value = 0.8477463511861221
$ws.Range("AL3").Value = 0.02460032727516325
$ws.Range("AM3").Value = 0.8406197842886998
$ws.Range("AN3").Value = 0.8068570687166512
$ws.Range("AO3").Value = 0.8612183458760958
$ws.Range("AP3").Value = 0.84881603261981
$ws.Range("AQ3").Value = 0.8812205244293529

# Row 4 values
$ws.Range("B4").Value = 0.8407444217568611
$ws.Range("C4").Value = 0.04745550508411889
$ws.Range("D4").Value = 0.8364762937870311
$ws.Range("E4").Value = 0.7874276261373035
$ws.Range("F4").Value = 0.8712616943516773
$ws.Range("G4").Value = 0.7946896811922279
$ws.Range("H4").Value = 0.9138668133160659
$ws.Range("I4").Value = 0.8452030091656443
$ws.Range("J4").Value = 0.02370525090879799
$ws.Range("K4").Value = 0.8710029809864569
$ws.Range("L4").Value = 0.8085731311537764
$ws.Range("M4").Value = 0.8595895144282241
$ws.Range("N4").Value = 0.8262635345247462
$ws.Range("O4").Value = 0.860585884735018
$ws.Range("P4").Value = 0.8413400401085502
$ws.Range("Q4").Value = 0.0359302007712412
$ws.Range("R4").Value = 0.8508379747420525
$ws.Range("S4").Value = 0.8171543641252531
$ws.Range("T4").Value = 0.8386922641591573
$ws.Range("U4").Value = 0.7969867267454821
$ws.Range("V4").Value = 0.9030288707708062
$ws.Range("W4").Value = 0.8296821578910329
$ws.Range("X4").Value = 0.0386187623816063
$ws.Range("Y4").Value = 0.8076050862165084
$ws.Range("Z4").Value = 0.7744879124305382
$ws.Range("AA4").Value = 0.8383175867046836
$ws.Range("AB4").Value = 0.8366758545768817
$ws.Range("AC4").Value = 0.8913243495265526
$ws.Range("AD4").Value = 0.8366381363445029
$ws.Range("AE4").Value = 0.04098591940013351
$ws.Range("AF4").Value = 0.8079261164367548
$ws.Range("AG4").Value = 0.7744879124305382
$ws.Range("AH4").Value = 0.8605254195930829
$ws.Range("AI4").Value = 0.8489268837355859
$ws.Range("AJ4").Value = 0.8913243495265526
$ws.Range("AK4").Value = 0.8389726415160503
$ws.Range("AL4").Value = 0.04910721072039775
$ws.Range("AM4").Value = 0.7971865586777432
$ws.Range("AN4").Value = 0.7744879124305382
$ws.Range("AO4").Value = 0.8605254195930829
$ws.Range("AP4").Value = 0.8489268837355859
$ws.Range("AQ4").Value = 0.9137364331433011

# Row 5 values
$ws.Range("B5").Value = 0.8130234233276126
$ws.Range("C5").Value = 0.03054598930282494
$ws.Range("D5").Value = 0.8497879528464635
$ws.Range("E5").Value = 0.7646768944959347
$ws.Range("F5").Value = 0.839085419730581
$ws.Range("G5").Value = 0.7959801850819846
$ws.Range("H5").Value = 0.8155866644830991
$ws.Range("I5").Value = 0.8459302328161922
$ws.Range("J5").Value = 0.03677573643999142
$ws.Range("K5").Value = 0.8503147305288963
$ws.Range("L5").Value = 0.7762299886611451
$ws.Range("M5").Value = 0.8812632802130176
$ws.Range("N5").Value = 0.8512515063651203
$ws.Range("O5").Value = 0.8705916583127821
$ws.Range("W5").Value = 0.8325940159128592
$ws.Range("X5").Value = 0.03501595222345227
$ws.Range("Y5").Value = 0.8694628253418847
$ws.Range("Z5").Value = 0.7803355461533162
$ws.Range("AA5").Value = 0.873337167106833
$ws.Range("AB5").Value = 0.8144893951345564
$ws.Range("AC5").Value = 0.8253451458277056
$ws.Range("AD5").Value = 0.8450343347222546
$ws.Range("AE5").Value = 0.02895801930486274
$ws.Range("AF5").Value = 0.8582585479137204
$ws.Range("AG5").Value = 0.799326116306677
$ws.Range("AH5").Value = 0.8828938441644972
$ws.Range("AI5").Value = 0.8268397633545902
$ws.Range("AJ5").Value = 0.8578534018717878

# Row 6 values
$ws.Range("B6").Value = 0.8209250479713978
$ws.Range("C6").Value = 0.04095143217120851
$ws.Range("D6").Value = 0.7834344745309959
$ws.Range("E6").Value = 0.7798217288331752
$ws.Range("F6").Value = 0.8186860234872239
$ws.Range("G6").Value = 0.8295761198987005
$ws.Range("H6").Value = 0.8931068931068931
$ws.Range("I6").Value = 0.7954769144201477
$ws.Range("J6").Value = 0.04217387143517815
$ws.Range("K6").Value = 0.7463044706979579
$ws.Range("L6").Value = 0.7779550767666218
$ws.Range("M6").Value = 0.8286864918522833
$ws.Range("N6").Value = 0.7647601028380433
$ws.Range("O6").Value = 0.8596784299458323
$ws.Range("P6").Value = 0.745153040415213
$ws.Range("Q6").Value = 0.0592943515002222
$ws.Range("R6").Value = 0.7030426908698092
$ws.Range("S6").Value = 0.6911444588863943
$ws.Range("T6").Value = 0.7734226269983444
$ws.Range("U6").Value = 0.7091207405920229
$ws.Range("V6").Value = 0.8490346847294947
$ws.Range("W6").Value = 0.8252049586505205
$ws.Range("X6").Value = 0.03789922721526721
$ws.Range("Y6").Value = 0.7682101306745358
$ws.Range("Z6").Value = 0.7994806042344241
$ws.Range("AA6").Value = 0.8605133589524849
$ws.Range("AB6").Value = 0.827685010707082
$ws.Range("AC6").Value = 0.8701356886840756
$ws.Range("AD6").Value = 0.8160231860916902
$ws.Range("AE6").Value = 0.02978961818739197
$ws.Range("AF6").Value = 0.7907374087235542
$ws.Range("AG6").Value = 0.827471572600466
$ws.Range("AH6").Value = 0.8601311972800005
$ws.Range("AI6").Value = 0.7757582242239992
$ws.Range("AJ6").Value = 0.8260175276304308
$ws.Range("AK6").Value = 0.8279993175251559
$ws.Range("AL6").Value = 0.03963272949201607
$ws.Range("AM6").Value = 0.793842062193126
$ws.Range("AN6").Value = 0.7755846745618579
$ws.Range("AO6").Value = 0.860415762225361
$ws.Range("AP6").Value = 0.8284075778666652
$ws.Range("AQ6").Value = 0.8817465107787688

# Row 7 values
$ws.Range("B7").Value = 0.8419954741168649
$ws.Range("C7").Value = 0.05187527408279901
$ws.Range("D7").Value = 0.8412209694404542
$ws.Range("E7").Value = 0.756633968433629
$ws.Range("F7").Value = 0.8701572507599673
$ws.Range("G7").Value = 0.828003280576886
$ws.Range("H7").Value = 0.9139619013733883
$ws.Range("I7").Value = 0.8640553429086986
$ws.Range("J7").Value = 0.03300848625555496
$ws.Range("K7").Value = 0.8270167032692771
$ws.Range("L7").Value = 0.8279620021555506
$ws.Range("M7").Value = 0.8808416380572578
$ws.Range("N7").Value = 0.8709677419354839
$ws.Range("O7").Value = 0.9134886291259225
$ws.Range("P7").Value = 0.8252465840462987
$ws.Range("Q7").Value = 0.04704288276541442
$ws.Range("R7").Value = 0.8000304100528066
$ws.Range("S7").Value = 0.7870787846806523
$ws.Range("T7").Value = 0.8303149589643105
$ws.Range("U7").Value = 0.7942110081329097
$ws.Range("V7").Value = 0.9145977584008144
$ws.Range("W7").Value = 0.8409894062204077
$ws.Range("X7").Value = 0.03610601058914542
$ws.Range("Y7").Value = 0.8295613789973716
$ws.Range("Z7").Value = 0.775164952190916
$ws.Range("AA7").Value = 0.8709061897102811
$ws.Range("AB7").Value = 0.8602180337112425
$ws.Range("AC7").Value = 0.8690964764922278
$ws.Range("AD7").Value = 0.8602931213053973
$ws.Range("AE7").Value = 0.03301800841502055
$ws.Range("AF7").Value = 0.8830129036765
$ws.Range("AG7").Value = 0.805988582364294
$ws.Range("AH7").Value = 0.8611559139784947
$ws.Range("AI7").Value = 0.8481485675034062
$ws.Range("AJ7").Value = 0.903159639004291
$ws.Range("AK7").Value = 0.8326997506786389
$ws.Range("AL7").Value = 0.02408657791928399
$ws.Range("AM7").Value = 0.8298341558658932
$ws.Range("AN7").Value = 0.795356024957543
$ws.Range("AO7").Value = 0.8289774741387643
$ws.Range("AP7").Value = 0.8385533994274889
$ws.Range("AQ7").Value = 0.8707776990035054
